$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 114, shifting existing rows 114-178 down to 115-179
$ws.Rows.Item(114).Insert()

# Populate the newly inserted row 114 with a new data record
$ws.Range("A114").Value = 10
$ws.Range("B114").Value = "Vega Modelo de Temuco"
$ws.Range("C114").Value = "La Araucanía"
$ws.Range("D114").Value = 44830
$ws.Range("E114").Value = 9
$ws.Range("F114").Value = "Fruta"
$ws.Range("G114").Value = 100104
$ws.Range("H114").Value = "Frutos de pepita"
$ws.Range("I114").Value = 100104001
$ws.Range("J114").Value = "Granada"
$ws.Range("K114").Value = "Wonderfull"
$ws.Range("L114").Value = "Primera"
$ws.Range("M114").Value = 300
$ws.Range("N114").Value = 14000
$ws.Range("O114").Value = 14000
$ws.Range("P114").Value = 14000
$ws.Range("Q114").Value = "`$/bandeja 10 kilos granel"
$ws.Range("R114").Value = "Provincia de Limarí"
$ws.Range("S114").Value = 1400
$ws.Range("T114").Value = 10
